$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the report date on the title page.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Date: 2024-09-12", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Date: 2024-09-13", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Overview of the nodes in the control flow" table.
#    - Drop the "Data Flow Task" row (its occurrence count survives onto the
#      "Execute SQL Task" row).
#    - Rename "Expression Task" to "Foreach Loop Container".
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$t1.Rows.Item(3).Delete()
$t1.Cell(2, 2).Range.Text = "1"
$t1.Cell(3, 1).Range.Text = "Foreach Loop Container"

# ---------------------------------------------------------------------------
# 3. "Overview of the nodes in the data flow" table.
#    Trim the trailing four rows and rewrite the surviving data rows with
#    the post-edit node/occurrence values.
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$t3.Rows.Item(11).Delete()
$t3.Rows.Item(10).Delete()
$t3.Rows.Item(9).Delete()
$t3.Rows.Item(8).Delete()

$t3.Cell(2, 1).Range.Text = "SSISODBCDst"
$t3.Cell(2, 2).Range.Text = "2"
$t3.Cell(3, 1).Range.Text = "DataSources"
$t3.Cell(3, 2).Range.Text = "2"
$t3.Cell(4, 1).Range.Text = "RowCount"
$t3.Cell(4, 2).Range.Text = "1"
$t3.Cell(5, 1).Range.Text = "ConditionalSplit"
$t3.Cell(5, 2).Range.Text = "1"
$t3.Cell(6, 1).Range.Text = "UnionAll"
$t3.Cell(6, 2).Range.Text = "1"
$t3.Cell(7, 1).Range.Text = "DataDestinations"
$t3.Cell(7, 2).Range.Text = "1"

# ---------------------------------------------------------------------------
# 4. "Overview of utilised source tables in the data flow" table.
#    Both source-table rows (Suppliers_Extract, Products_Extract) are gone.
# ---------------------------------------------------------------------------
$t4 = $d.Tables.Item(4)
$t4.Rows.Item(3).Delete()
$t4.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# 5. "Overview of utilised target tables in the data flow" table.
#    "Error_lines" becomes "Supp_Prod_output".
# ---------------------------------------------------------------------------
$t5 = $d.Tables.Item(5)
$t5.Cell(2, 1).Range.Text = "Supp_Prod_output"

# ---------------------------------------------------------------------------
# 6. Node / Column / Transformation table.
# ---------------------------------------------------------------------------
$t6 = $d.Tables.Item(6)
$t6.Cell(2, 1).Range.Text = "Merge and filter@Row Count errors"
$t6.Cell(3, 1).Range.Text = "Supp_Prod_output"
